# Add files via upload
# Replace the player roster table (A2:C19) with the new data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Luke Kennard", "SG", "Memphis Grizzlies"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Bam Adebayo", "C", "Miami Heat"),
    @("Wendell Carter Jr.", "PF,C", "Orlando Magic"),
    @("Nikola Jovic", "PF,C", "Miami Heat"),
    @("Toumani Camara", "SF,PF", "Portland Trail Blazers"),
    @("Jordan Clarkson", "SG,SF", "Utah Jazz"),
    @("Shai Gilgeous-Alexander", "PG,SG", "Oklahoma City Thunder"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Jose Alvarado", "PG", "New Orleans Pelicans"),
    @("Mike Conley", "PG", "Minnesota Timberwolves"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
